# Moving from 3.1.1 to 3.2.0.
#
# Insert an M2Doc "version mismatch" warning block right after the
# "Basic " run and before the "if" run on the first paragraph, mirroring
# the existing invalid-expression warning block further down in the
# document (same layout: 4 spaces, "<---", message, 4 spaces) but using
# an orange (FFA500) color instead of red (FF0000).

$d = $word.ActiveDocument

$marker = "<---"
$message = "M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0"

# Locate the "Basic " text at the very start of the document and collapse
# the found range to its end so the new block is inserted right after it
# (and before the "if" that currently follows it).
$rng = $d.Content
$found = $rng.Find.Execute("Basic ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Basic ' text to anchor the warning block"
}
$rng.Collapse(0)
$insertStart = $rng.Start

# Insert the whole block as plain text first (4 spaces, marker, message,
# 4 spaces) so nothing accidentally inherits stray character formatting.
$rng.InsertAfter("    " + $marker + $message + "    ")

# Now reach for just the "<---message" portion (not the surrounding
# spaces) through a fresh, precisely bounded Range and format only that
# part: orange color, 16pt, light gray highlight - matching the existing
# invalid-expression warning elsewhere in the document.
$colorStart = $insertStart + 4
$colorEnd = $colorStart + ($marker + $message).Length
$colorRng = $d.Range($colorStart, $colorEnd)
$colorRng.Font.Color = 42495
$colorRng.Font.Size = 16
$colorRng.Font.HighlightColorIndex = 16
